$p = $ppt.ActivePresentation

# --------------------------------------------------------------------
# 1) Table on slide 16 ("PLENARY- COMPLETE THE MISSING GAPS") switches
#    from the built-in "Table_0" style to another built-in table style.
# --------------------------------------------------------------------
$tableSlide = $p.Slides.Item(16)
for ($i = 1; $i -le $tableSlide.Shapes.Count; $i++) {
    $shp = $tableSlide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{02CE31C4-319C-4FD1-96D5-55786E6D55A5}")
    }
}

# --------------------------------------------------------------------
# 2) Theme colours: the deck's theme colour scheme (currently the
#    "Integral" palette) is changed to the plain "Office Theme"
#    palette (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink).
# --------------------------------------------------------------------
$officeThemeRGB = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

$slide1 = $p.Slides.Item(1)
$themeColors = $slide1.ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Item($i).RGB = $officeThemeRGB[$i - 1]
}
